$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.290.13'
$ws.Range('E2').Value = '  +2.33%  '
$ws.Range('D3').Value = '2.095.60'
$ws.Range('E3').Value = '  +4.41%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.61'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.24%  '
$ws.Range('E6').Value = '  +0.43%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '53.68'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +19.57%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '61.74'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.72%  '
$ws.Range('E10').Value = '  +1.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0742'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.05%  '
$ws.Range('E12').Value = '  +7.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.31'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.09%  '
$ws.Range('D14').Value = '2.403.69'
$ws.Range('E14').Value = '  +4.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.839'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.13%  '
$ws.Range('D16').Value = '2.093.53'
$ws.Range('E16').Value = '  +4.20%  '
$ws.Range('E17').Value = '  +5.16%  '
$ws.Range('D18').Value = '37.217.73'
$ws.Range('E18').Value = '  +2.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.64'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.57'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +13.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '241.15'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.21'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +6.87%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('E25').Value = '  +2.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '171.75'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.30'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.67'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.28%  '
$ws.Range('E29').Value = '  +4.34%  '
$ws.Range('E30').Value = '  +1.66%  '
$ws.Range('E31').Value = '  +26.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '22.14'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.51'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.98%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0613'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0906'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +12.87%  '
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.30'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.31%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.12'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.47%  '
$ws.Range('B39').Value = 'WEMIXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.85'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.34'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.52'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +15.97%  '
$ws.Range('E42').Value = '  +3.98%  '
$ws.Range('E43').Value = '  +5.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '98.85'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0922'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +12.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.79'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.13'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +100.81%  '
$ws.Range('D48').Value = '1.317.76'
$ws.Range('E48').Value = '  +0.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.96'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.04'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +14.58%  '
$ws.Range('D51').Value = '2.292.32'
$ws.Range('E51').Value = '  +4.11%  '
